# Sort the comma-separated "Recorded By" names in column G ("Recorded By")
# of the "Session Analysis Results" sheet using an ordinal (ASCII, case
# sensitive, uppercase-before-lowercase) alphabetical sort - matching the
# behavior of e.g. Python's sorted() / C's strcmp, which the default
# PowerShell string comparison operators in this runtime do not reproduce
# (they are case-insensitive).

function Sort-Ordinal($items) {
    $arr = @($items)
    $n = $arr.Count
    for ($i = 1; $i -lt $n; $i++) {
        $key = $arr[$i]
        $j = $i - 1
        while ($j -ge 0) {
            $s1 = $arr[$j]
            $s2 = $key
            $len1 = $s1.Length
            $len2 = $s2.Length
            $minLen = [Math]::Min($len1, $len2)
            $cmp = 0
            for ($k = 0; $k -lt $minLen; $k++) {
                $c1 = [int][char]$s1[$k]
                $c2 = [int][char]$s2[$k]
                if ($c1 -ne $c2) {
                    $cmp = $c1 - $c2
                    break
                }
            }
            if ($cmp -eq 0) {
                $cmp = $len1 - $len2
            }
            if ($cmp -le 0) {
                break
            }
            $arr[$j+1] = $arr[$j]
            $j = $j - 1
        }
        $arr[$j+1] = $key
    }
    return $arr
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $sortedParts = Sort-Ordinal $parts
            $sortedText = [string]::Join(", ", $sortedParts)
            if ($sortedText -ne $val) {
                $cell.Value2 = $sortedText
            }
        }
    }
}
